$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Remove the "Statistical Inference" / "Maximum likelihood" rows (old rows 12-13) ---
# This shifts every row below up by two, and Excel automatically keeps the
# F17 formula (now K15) and shared-string table consistent.
$ws.Rows("12:13").Delete()

# --- Trim / update the remaining "Pages" (column F) numbers ---
$ws.Range("F3").Value  = 21
$ws.Range("F5").Value  = 2.05
$ws.Range("F6").Value  = 2
$ws.Range("F7").Value  = 0.1
$ws.Range("F8").Value  = 1.3
$ws.Range("F9").Value  = 0.45
$ws.Range("F10").Value = 0.45
$ws.Range("F11").Value = 0.1
$ws.Range("F12").Value = 0.4
$ws.Range("F13").Value = 3.2
$ws.Range("F17").Value = 0.7
$ws.Range("F19").Value = 0.4
$ws.Range("F20").Value = 1.1
$ws.Range("F22").Value = 0.7
$ws.Range("F23").Value = 0.95
$ws.Range("F25").Value = 0.6
$ws.Range("F26").Value = 2.5
$ws.Range("F29").Value = 1

# --- Add the running-total check formula next to the grand total ---
$ws.Range("G3").Formula = "=SUM(F4,F5,F6,F13,F20,F23,F26,F27,F28,F29,F30)"

# --- Update the saved view: select F6 (also clears the old scrolled-down view) ---
$ws.Range("F6").Select() | Out-Null

Write-Output "done"
